# Reporting.xlsx – Phase 2 reflection edits
# (Bugfix: Partial csv code implemented in budgets view removed. Phase 2 reflection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename three existing "Create Sort Function for ... Class" entries to
#    "... Controller" (sort functions actually live on the controllers).
# ---------------------------------------------------------------------------
$ws.Range("E43").Value = "Create Sort Function for Transaction Controller"
$ws.Range("E44").Value = "Create Sort Function for Tag Controller"
$ws.Range("E45").Value = "Create Sort Function for Merchant Controller"

# ---------------------------------------------------------------------------
# 2. Row 46 was missing its Manhours entry (0.5) and used the wrong
#    (non-centered) style for column C - bring it in line with the other
#    rows in this block.
# ---------------------------------------------------------------------------
$ws.Range("B46").Value = 0.5
$ws.Range("B46").NumberFormat = $ws.Range("B44").NumberFormat()
$ws.Range("C46").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 3. New row 47 - Budget Controller sort function.
# ---------------------------------------------------------------------------
$ws.Range("A47").Value = 43508
$ws.Range("A47").NumberFormat = $ws.Range("A46").NumberFormat()

$ws.Range("B47").Value = 0.5
$ws.Range("B47").NumberFormat = $ws.Range("B44").NumberFormat()

$ws.Range("C47").Value = $ws.Range("C46").Value()
$ws.Range("C47").HorizontalAlignment = -4108

$ws.Range("D47").Value = 2

$ws.Range("E47").Value = "Create Sort Function for Budget Controller"
$ws.Range("F47").Value = "Once you’ve done one...are sorts always in controller?"

# ---------------------------------------------------------------------------
# 4. New row 48 - Phase 3 reflection / CSV import-export note.
# ---------------------------------------------------------------------------
$ws.Range("A48").Value = 43508
$ws.Range("A48").NumberFormat = $ws.Range("A46").NumberFormat()

$ws.Range("B48").Value = 0.5
$ws.Range("B48").NumberFormat = $ws.Range("B44").NumberFormat()

$ws.Range("C48").Value = "Reflection"
$ws.Range("C48").HorizontalAlignment = -4108

$ws.Range("D48").Value = 2

$ws.Range("E48").Value = "Review views for phase 3"
$ws.Range("F48").Value = "Put CSV import/export on hold until CSS done, product is stable."

# ---------------------------------------------------------------------------
# 5. Update the selection / scroll position to match the author's final view.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 1
$ws.Range("E47").Select()
